$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New numeric value added next to the existing header/data area (row 2, column E)
$ws.Range("E2").Value = 12

# Append 18 more data rows (38-55) that repeat the same user record already
# present in rows 11-37 (userId=user112, firstName=giri, lastName=t, password=passwor)
for ($r = 38; $r -le 55; $r++) {
  $ws.Cells.Item($r, 1).Value = "user112"
  $ws.Cells.Item($r, 2).Value = "giri"
  $ws.Cells.Item($r, 3).Value = "t"
  $ws.Cells.Item($r, 4).Value = "passwor"
}

# Move the active selection to the newly-added cell
$ws.Range("E2").Select() | Out-Null
